# Sentencing-date corrections for the Fine Only Judgment Entry.
# The change-of-plea hearing, the fines/costs due date, and the license
# suspension start date all move from June 12 to June 13, 2022; the
# community-service completion deadline moves from August 11 to August 12.
$d = $word.ActiveDocument

# "...for a change of plea on June 12, 2022." (hearing date - its own run)
$d.Content.Find.Execute(" on June 12, 2022.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " on June 13, 2022.", 2)

# Remaining bare "June 12, 2022" occurrences each live alone inside their own
# run (the bold "in full by <date>" run, and the "license is suspended from
# <date>" run), so one case-sensitive ReplaceAll cleanly updates both without
# touching neighboring runs/formatting.
$d.Content.Find.Execute("June 12, 2022", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "June 13, 2022", 2)

# "...Office of Community Control on or before August 11, 2022." (community service deadline)
$d.Content.Find.Execute("August 11, 2022", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "August 12, 2022", 2)
